# Adds season-record columns (Wins / Losses / Ties) to the sheet.
# This mirrors the upstream change where the scraper was fixed to also
# pull down the team's season record, not just individual stats.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new column headers in AD1:AF1 ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/centered/bordered header style (same as A1:AC1)
# by copying the formatting from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-49): season record values for every player row ---
$lastRow = 49
$wins = $ws.Range("AD2:AD" + $lastRow)
$losses = $ws.Range("AE2:AE" + $lastRow)
$ties = $ws.Range("AF2:AF" + $lastRow)

$wins.Value = 74
$losses.Value = 87
$ties.Value = 0
